$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 0.019999999552965164
$ws.Range("E4").Value = 0.029999999329447746
$ws.Range("E5").Value = 0.039999999105930328
$ws.Range("E6").Value = 0.05000000074505806
$ws.Range("E7").Value = 0.059999998658895493
$ws.Range("E8").Value = 0.070000000298023224
$ws.Range("E9").Value = 0.079999998211860657
$ws.Range("E10").Value = 0.090000003576278687
$ws.Range("E11").Value = 0.10000000149011612
$ws.Range("E12").Value = 0.10999999940395355
$ws.Range("E13").Value = 0.11999999731779099
